$d = $word.ActiveDocument

$replacements = @(
    @{ old = "2025-10-14 Tuesday"; new = "2025-10-15 Wednesday" },
    @{ old = "987×4=3948"; new = "111×5=555" },
    @{ old = "232×7=1624"; new = "841×8=6728" },
    @{ old = "645×7=4515"; new = "680×5=3400" },
    @{ old = "498×2=996"; new = "274×7=1918" },
    @{ old = "306×5=1530"; new = "942×6=5652" },
    @{ old = "831×6=4986"; new = "149×8=1192" },
    @{ old = "386×6=2316"; new = "770×8=6160" },
    @{ old = "350×4=1400"; new = "619×5=3095" },
    @{ old = "872×3=2616"; new = "872×9=7848" },
    @{ old = "742×8=5936"; new = "579×6=3474" },
    @{ old = "892×3=2676"; new = "855×3=2565" },
    @{ old = "891×3=2673"; new = "936×3=2808" },
    @{ old = "676×5=3380"; new = "732×8=5856" },
    @{ old = "121×5=605"; new = "696×6=4176" },
    @{ old = "832×3=2496"; new = "612×8=4896" },
    @{ old = "566×9=5094"; new = "969×2=1938" },
    @{ old = "894×2=1788"; new = "662×3=1986" },
    @{ old = "857×7=5999"; new = "860×5=4300" },
    @{ old = "121×6=726"; new = "698×6=4188" },
    @{ old = "784×5=3920"; new = "672×9=6048" },
    @{ old = "634×5=3170"; new = "925×5=4625" },
    @{ old = "998×4=3992"; new = "178×2=356" },
    @{ old = "922×7=6454"; new = "335×8=2680" },
    @{ old = "233×8=1864"; new = "721×5=3605" },
    @{ old = "900×4=3600"; new = "155×4=620" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
